# "major accuracy check update"
#
# The harvestDate column (A2:A37) currently reads "9.13.19" and needs to be
# corrected to the zero-padded "09.13.19" everywhere it appears. The column
# is a shared-string text value (not a real date), so we have to make sure
# Excel's automatic "this looks like a date" parser doesn't silently turn it
# into a date serial number when we write the new text back in.
#
# Trick: temporarily force the range to Text format ("@") while we assign the
# corrected string, then restore the original (General) number format. That
# keeps the cells as plain text with their original (default) style, exactly
# like before the edit, while the stored/shared string itself is corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("A2:A37")

$originalFormat = $dateRange.NumberFormat

$dateRange.NumberFormat = "@"
$dateRange.Value = "09.13.19"
$dateRange.NumberFormat = $originalFormat

# Reflect where the editor's cursor ended up after making the change.
$ws.Range("A38").Select()
